$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# Row 4: Qty Storage Nodes changes from 6 to 4
$ws.Range("W4").Value = 4

# Row 5: update existing values
$ws.Range("D5").Value = "Test-Automation-Indirect"
$ws.Range("J5").Value = "Indirect"

# Row 5: fill in previously empty Partner Sales Model related columns
$ws.Range("L5").Value = "Reseller through NetApp"
$ws.Range("M5").Value = "Reseller"
$ws.Range("N5").Value = "Own Use"
$ws.Range("O5").Value = "Test"
$ws.Range("P5").Value = "Bhuvan Testing"
$ws.Range("Q5").Value = "Bhuvan Testing"
